# Apply latest crypto price/volume updates (GitHub Actions refresh).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'28.639.49"
$ws.Range("E2").Value = "  +1.02%  "
$ws.Range("D3").Value = "'1.561.59"
$ws.Range("E3").Value = "  -0.57%  "
$ws.Range("D4").Value = "'0.998"
$ws.Range("E4").Value = "  -0.53%  "
$ws.Range("D5").Value = "'210.29"
$ws.Range("E5").Value = "  -0.55%  "
$ws.Range("E6").Value = "  -1.02%  "
$ws.Range("D7").Value = "'0.998"
$ws.Range("E7").Value = "  -0.52%  "
$ws.Range("D8").Value = "'25.13"
$ws.Range("E8").Value = "  +5.80%  "
$ws.Range("E9").Value = "  -0.75%  "
$ws.Range("E10").Value = "  -0.40%  "
$ws.Range("D11").Value = "'0.0894"
$ws.Range("E11").Value = "  -0.08%  "
$ws.Range("D12").Value = "'1.784.73"
$ws.Range("D13").Value = "'1.558.25"
$ws.Range("E13").Value = "  -0.78%  "
$ws.Range("D14").Value = "'28.642.49"
$ws.Range("E14").Value = "  +0.97%  "
$ws.Range("E15").Value = "  +0.13%  "
$ws.Range("E16").Value = "  -1.11%  "
$ws.Range("D17").Value = "'61.20"
$ws.Range("E17").Value = "  -0.31%  "
$ws.Range("D18").Value = "'228.06"
$ws.Range("E18").Value = "  +0.08%  "
$ws.Range("E19").Value = "  -0.69%  "
$ws.Range("E20").Value = "  -1.13%  "
$ws.Range("D21").Value = "'0.998"
$ws.Range("E21").Value = "  -0.46%  "
$ws.Range("D22").Value = "'3.90"
$ws.Range("E23").Value = "  +0.34%  "
$ws.Range("E24").Value = "  +0.98%  "
$ws.Range("D25").Value = "'151.19"
$ws.Range("D26").Value = "'14.74"
$ws.Range("E26").Value = "  -1.19%  "
$ws.Range("E27").Value = "  +0.23%  "
$ws.Range("D28").Value = "'0.999"
$ws.Range("E28").Value = "  -0.46%  "
$ws.Range("E29").Value = "  -1.88%  "
$ws.Range("D30").Value = "'0.0460"
$ws.Range("E30").Value = "  -4.18%  "
$ws.Range("E31").Value = "  -2.60%  "
$ws.Range("E32").Value = "  -0.16%  "
$ws.Range("D33").Value = "'1.389.85"
$ws.Range("E33").Value = "  +0.72%  "
$ws.Range("E34").Value = "  -3.88%  "
$ws.Range("E35").Value = "  -4.53%  "
$ws.Range("E36").Value = "  -1.48%  "
$ws.Range("E37").Value = "  +1.55%  "
$ws.Range("E38").Value = "  -2.33%  "
$ws.Range("E40").Value = "  +2.37%  "
$ws.Range("E41").Value = "  -0.82%  "
$ws.Range("E42").Value = "  -0.43%  "
$ws.Range("D43").Value = "'0.769"
$ws.Range("E43").Value = "  -1.92%  "
$ws.Range("D44").Value = "'0.0458"
$ws.Range("E44").Value = "  -2.61%  "
$ws.Range("D45").Value = "'63.98"
$ws.Range("E45").Value = "  +2.71%  "
$ws.Range("D46").Value = "'5.22"
$ws.Range("E46").Value = "  -2.18%  "
$ws.Range("D47").Value = "'1.696.94"
$ws.Range("E47").Value = "  -0.68%  "
$ws.Range("E48").Value = "  -5.31%  "
$ws.Range("D49").Value = "'84.99"
$ws.Range("E49").Value = "  -0.47%  "
$ws.Range("D50").Value = "'43.22"
$ws.Range("E50").Value = "  +7.05%  "
$ws.Range("E51").Value = "  -0.65%  "
